$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give J1 the same format as the rest of the header row (centered style)
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

# New column J: "estado" header + status values per row
$ws.Range("J1").Value = "estado"
$ws.Range("J2").Value = "OPCC"
$ws.Range("J3").Value = "OPCC"
$ws.Range("J4").Value = "CBS"
$ws.Range("J5").Value = "CBS"
$ws.Range("J6").Value = "CBS"
$ws.Range("J7").Value = "CBS"
$ws.Range("J8").Value = "TRUCK"
$ws.Range("J9").Value = "TRUCK"
$ws.Range("J10").Value = "INVENT"

$excel.CutCopyMode = $false
$ws.Range("H17").Select()
